$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-07 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-08 Sunday", 2)

$d.Content.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷6=7, 4", 2)
$d.Content.Find.Execute("85÷9=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "67÷2=33, 1", 2)
$d.Content.Find.Execute("41÷6=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "52÷2=26, 0", 2)
$d.Content.Find.Execute("14÷2=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷2=36, 0", 2)
$d.Content.Find.Execute("23÷2=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "96÷9=10, 6", 2)

$d.Content.Find.Execute("36÷8=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=7, 6", 2)
$d.Content.Find.Execute("63÷8=7, 7", $true, $false, $false, $false, $false, $true, 1, $false, "51÷8=6, 3", 2)
$d.Content.Find.Execute("59÷7=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "66÷7=9, 3", 2)
$d.Content.Find.Execute("63÷4=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "30÷4=7, 2", 2)
$d.Content.Find.Execute("59÷4=14, 3", $true, $false, $false, $false, $false, $true, 1, $false, "70÷9=7, 7", 2)

$d.Content.Find.Execute("22÷9=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "43÷4=10, 3", 2)
$d.Content.Find.Execute("65÷8=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "26÷2=13, 0", 2)
$d.Content.Find.Execute("89÷5=17, 4", $true, $false, $false, $false, $false, $true, 1, $false, "40÷2=20, 0", 2)
$d.Content.Find.Execute("78÷3=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷4=20, 2", 2)
$d.Content.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "23÷9=2, 5", 2)

$d.Content.Find.Execute("37÷3=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "96÷9=10, 6", 2)
$d.Content.Find.Execute("87÷4=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "56÷6=9, 2", 2)
$d.Content.Find.Execute("25÷3=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷5=12, 0", 2)
$d.Content.Find.Execute("50÷4=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=17, 0", 2)
$d.Content.Find.Execute("91÷8=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "34÷8=4, 2", 2)

$d.Content.Find.Execute("30÷3=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "41÷7=5, 6", 2)
$d.Content.Find.Execute("20÷2=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "71÷5=14, 1", 2)
$d.Content.Find.Execute("95÷9=10, 5", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=22, 1", 2)
$d.Content.Find.Execute("39÷5=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "97÷6=16, 1", 2)
$d.Content.Find.Execute("76÷7=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "96÷6=16, 0", 2)
